# Updates crypto price/volume data to match the latest scrape.
# Column D ("Price") and E ("Volume(1h)") values are stored as literal text
# (e.g. "53.849.44", "  -4.54%  ") rather than numbers, matching the source
# workbook which uses inline/shared strings for every data cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the default (unformatted) style so we can restore it after
# temporarily forcing Text format on cells whose new value looks numeric
# (Excel would otherwise auto-convert "487.93" etc. into a real number).
$defaultStyle = $ws.Range("D5").Style

$ws.Range("D2").Value = "53.849.44"
$ws.Range("E2").Value = "  -4.54%  "
$ws.Range("D3").Value = "2.235.48"
$ws.Range("E3").Value = "  -6.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.93"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.66"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").Value = "2.247.82"
$ws.Range("E9").Value = "  -5.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0917"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -6.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.72"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "2.630.97"
$ws.Range("E14").Value = "  -6.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.16"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "53.769.94"
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "2.246.57"
$ws.Range("E18").Value = "  -9.03%  "
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("E20").Value = "  -4.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "299.44"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.72"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.364"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.06"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.03"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "0.0₃0690"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.75"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +6.44%  "
$ws.Range("E39").Value = "  -5.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.64"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.80"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -6.30%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.68"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0878"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.536"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "236.74"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.14"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -4.41%  "
